$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 1 header cell updates ---
$ws.Range("E1").Value = "symbol"
$ws.Range("F1").Value = "name"
$ws.Range("G1").Value = "sourcespecies "
$ws.Range("H1").Value = "editedname"

# --- Row 2 data cell updates ---
$ws.Range("B2").Value = "check that user can search by edited GF>> Name"
$ws.Range("C2").Value = "1 mtihvidtst ifmrkayyen mvtipevvde irdensqfyf sllnlrveea snrnvekvir  61 vakktgdihk lsntdiklia kaldikerge dvilvtddys iqnvamslgl kvdnivqpki 121 skrfrwvkvc rgcgrsvdgd icpvcgseam ikkvrr"
$ws.Range("E2").Value = "GF_name_edit"
$ws.Range("F2").Value = "editName_GF"
$ws.Range("H2").Value = "name_GF_updated"

# --- Row height ---
$ws.Rows.Item(2).RowHeight = 45

# --- Column widths ---
$ws.Columns.Item(5).ColumnWidth = 18.6
$ws.Range("H1:I1").ColumnWidth = 20.5

# --- Selection / view state ---
[void]$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("F1").Select()

# --- Workbook window geometry (best effort) ---
$excel.ActiveWindow.Top = 1185
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Width = 12690
$excel.ActiveWindow.Height = 6555
